$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell C1 - same style as B1 (bold, bordered, centered header)
$ws.Range("C1").Value = 2
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122) # xlPasteFormats

# Data values for the new "t+3" column C
$values = @(
    -5.022459607464143,
    -1.15825194988682,
    -0.07204406301364299,
    -0.4144953840754857,
    0.01451842867919532,
    0.1067097157949464,
    0.1344667699115433,
    0.02767546902356237
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
